# data-mining: export excel as pdf
# Re-lay the little frequency-table block: it used to start at D3 (header
# row) / D4 (first data row); move it one column left and two rows up so it
# starts at C1 (header) / C2 (first data row). Formulas are re-typed against
# their new neighbours instead of being shared-formula copies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. wipe the old D3:H10 block -----------------------------------------
$ws.Range("D3:H10").ClearContents()

# --- 2. headers, now on row 1, columns C:G ---------------------------------
$ws.Range("C1").Value = "Значение"
$ws.Range("D1").Value = "Частота"
$ws.Range("E1").Value = "Частность"
$ws.Range("F1").Value = "Накопленная частота"
$ws.Range("G1").Value = "Накопленная частость"

# --- 3. data rows 2..8, columns C:G -----------------------------------------
$values = @(6, 7, 8, 9, 10, 11, 12)
for ($i = 0; $i -lt $values.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 3).Value = $values[$i]
    $ws.Cells.Item($r, 4).Formula = "=COUNTIF(`$A`$1:`$A`$47,C$r)"
    $ws.Cells.Item($r, 5).Formula = "=D$r/SUM(`$D`$1:`$D`$10)"
    $ws.Cells.Item($r, 6).Formula = "=SUM(D`$2:D$r)"
    $ws.Cells.Item($r, 7).Formula = "=SUM(E`$2:E$r)"
}

# --- 4. defined name "Extract" now covers only the last 10 rows ------------
$wb.Names.Item("Sheet1!Extract").RefersTo = '=Sheet1!$C$9:$C$18'

# --- 5. charts: series formulas follow the moved "cumulative" columns ------
$charts = $ws.ChartObjects()

$chart1 = $charts.Item(1).Chart
$chart1.SeriesCollection().Item(1).Formula = '=SERIES(,,Sheet1!$F$2:$F$8,1)'

$chart2 = $charts.Item(2).Chart
$chart2.SeriesCollection().Item(1).Formula = '=SERIES(,,Sheet1!$G$2:$G$8,1)'

# --- 6. charts: re-position/resize to sit below the (now higher) table -----
$c1 = $charts.Item(1)
$c1.Left = 141.99996062992125
$c1.Top = 686.4999212598425
$c1.Width = 346.30570343257875
$c1.Height = 221.2500787401575

$c2 = $charts.Item(2)
$c2.Left = 143.375
$c2.Top = 375.75
$c2.Width = 345.9306640625
$c2.Height = 224.0

$c3 = $charts.Item(3)
$c3.Left = 146.125
$c3.Top = 130.5
$c3.Width = 341.1806640625
$c3.Height = 222.25

# --- 7. selection, matching where the user ended up working ----------------
$ws.Range("I43").Select()
